# Junction_Flooding_247.xlsx edit:
#  - Round row 5 (B5:AH5) values to 2 decimal places ("custom accuracy")
#  - Remove row 6 (reduce the sampled dataset)
#  - Dimension will automatically shrink from A1:AH6 to A1:AH5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the data values in row 5 to 2 decimal places.
$row5Values = @{
    "B5"  = 16.19
    "C5"  = 12.1
    "D5"  = 1.07
    "E5"  = 35.51
    "F5"  = 28.82
    "G5"  = 12.63
    "H5"  = 50.96
    "I5"  = 19.73
    "J5"  = 8.89
    "K5"  = 12.74
    "L5"  = 14.26
    "M5"  = 15.18
    "N5"  = 4.11
    "O5"  = 12.78
    "P5"  = 18.09
    "Q5"  = 10.89
    "R5"  = 0.69
    "S5"  = 0.57
    "T5"  = 187.25
    "U5"  = 35.79
    "V5"  = 11.8
    "W5"  = 23.94
    "X5"  = 12.58
    "Y5"  = 1.7
    "Z5"  = 25.08
    "AA5" = 10.42
    "AB5" = 9.27
    "AC5" = 10.9
    "AD5" = 15.01
    "AE5" = 0.52
    "AF5" = 46.49
    "AG5" = 6.57
    "AH5" = 14.75
}

foreach ($addr in $row5Values.Keys) {
    $ws.Range($addr).Value = $row5Values[$addr]
}

# Remove row 6 entirely (data trimmed to fewer sampled rows).
$ws.Rows.Item(6).Delete()
